$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.062.67'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '1.834.35'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('D4').Value = '''0.9995'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''244.42'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').Value = '''0.6349'
$ws.Range('E6').Value = '  +2.06%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '''0.07560'
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('D9').Value = '''0.2950'
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('D10').Value = '''22.97'
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('D11').Value = '''0.07741'
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('D12').Value = '1.836.22'
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('D13').Value = '''5.014'
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '''0.6723'
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('D15').Value = '''83.27'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').Value = '''0.000009728'
$ws.Range('E16').Value = '  +7.22%  '
$ws.Range('D17').Value = '''6.099'
$ws.Range('E17').Value = '  +1.75%  '
$ws.Range('D18').Value = '29.087.81'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').Value = '''12.61'
$ws.Range('E19').Value = '  +2.36%  '
$ws.Range('D20').Value = '''226.52'
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('D21').Value = '''0.9995'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').Value = '''7.206'
$ws.Range('E22').Value = '  +0.25%  '
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '''160.57'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('D25').Value = '''0.1405'
$ws.Range('E25').Value = '  +3.21%  '
$ws.Range('D26').Value = '''8.560'
$ws.Range('E26').Value = '  +1.80%  '
$ws.Range('D27').Value = '''17.94'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').Value = '''1.500'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').Value = '''4.129'
$ws.Range('E29').Value = '  +2.09%  '
$ws.Range('D30').Value = '''4.091'
$ws.Range('E30').Value = '  +1.50%  '
$ws.Range('D31').Value = '''1.204'
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('D32').Value = '''0.05399'
$ws.Range('E32').Value = '  +3.54%  '
$ws.Range('D33').Value = '''1.868'
$ws.Range('E33').Value = '  +1.97%  '
$ws.Range('D34').Value = '''0.7464'
$ws.Range('E34').Value = '  +1.94%  '
$ws.Range('D35').Value = '''1.145'
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('D36').Value = '''2.657'
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('D37').Value = '1.246.15'
$ws.Range('E37').Value = '  -2.15%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '''2.759'
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.01790'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').Value = '''6.658'
$ws.Range('E40').Value = '  +5.32%  '
$ws.Range('D41').Value = '''0.9075'
$ws.Range('E41').Value = '  +1.72%  '
$ws.Range('D42').Value = '''0.9999'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '''101.98'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = '1.985.11'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').Value = '''65.00'
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('E46').Value = '  +2.79%  '
$ws.Range('D47').Value = '''0.5114'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '''0.4099'
$ws.Range('E48').Value = '  +3.48%  '
$ws.Range('D49').Value = '''9.057'
$ws.Range('E49').Value = '  +2.30%  '
$ws.Range('D50').Value = '''6.786'
$ws.Range('E50').Value = '  +2.19%  '
$ws.Range('D51').Value = '''1.650'
$ws.Range('E51').Value = '  -0.94%  '
